$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set target cells to Text format first so numeric-looking strings
# (e.g. "319.53", "4.76%") are stored as literal text, matching the
# original inlineStr cells rather than being parsed into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "319.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.76%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.61%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.320"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.83%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08047"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.62%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.597"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.70%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "27.38%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.648"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1271"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.55%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1959"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.77%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09522"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.61%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04522"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "9.11%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1046"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.06%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001299"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04208"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.77%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005914"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.84%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.347"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.66%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.476"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "5.62%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3481"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.40%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.211"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.30%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.41%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001296"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.68%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004351"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.00%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.97%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003541"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.24%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02705"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "5.01%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05920"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "10.56%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "93.50%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008045"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.83%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1466"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.52%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007526"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.60%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007954"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.33%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3212"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.10%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007010"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.92%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.94%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05582"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-8.50%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004002"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.77%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.94%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.94%"
